# Automated "Disponibilidad" refresh run (Actualizar 02-04-2021 18-35-09).
# The monitoring script re-checks each endpoint and appends a fresh batch
# of 14 rows (same Nombre/URL/Disponibilidad cycle) stamped with the
# current run's timestamp; it also nudges the immediately preceding
# batch's timestamp to the value recomputed during this same run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The previous batch (rows 58:71) gets its "Fecha" timestamp refreshed
# to the value recomputed by this run.
$ws.Range("D58:D71").Value = 44231.76464553241

# Append the next 14-row batch by duplicating the last block (rows
# 58:71 -> 72:85); this carries over the Nombre/URL/Disponibilidad
# values, number formats and styles untouched.
$src = $ws.Range("A58:D71")
$dst = $ws.Range("A72")
$src.Copy($dst)

# Stamp the new batch with this run's timestamp.
$ws.Range("D72:D85").Value = 44231.77435050038

# Re-create the hyperlinks on the new batch's URL column (Copy/Paste
# does not carry hyperlinks across), matching the same target cycle as
# every previous batch. The MapStore link (9th row of the cycle) keeps
# its "/" fragment sub-address.
$ws.Hyperlinks.Add($ws.Range("B72"), "https://www.dataintelligence-group.com/")
$ws.Hyperlinks.Add($ws.Range("B73"), "https://serviciodashboard.azurewebsites.net/")
$ws.Hyperlinks.Add($ws.Range("B74"), "https://powerbi.microsoft.com/es-es/")
$ws.Hyperlinks.Add($ws.Range("B75"), "https://www.dropbox.com/")
$ws.Hyperlinks.Add($ws.Range("B76"), "https://dataintelligence.store/")
$ws.Hyperlinks.Add($ws.Range("B77"), "https://app-data-i.users.earthengine.app/")
$ws.Hyperlinks.Add($ws.Range("B78"), "https://odooutil.azurewebsites.net/")
$ws.Hyperlinks.Add($ws.Range("B79"), "https://filtradordashboard.azurewebsites.net/")
$ws.Hyperlinks.Add($ws.Range("B80"), "https://ide.dataintelligence-group.com/mapstore/", "/")
$ws.Hyperlinks.Add($ws.Range("B81"), "https://ide.dataintelligence-group.com/geoserver/web/?0")
$ws.Hyperlinks.Add($ws.Range("B82"), "https://ide.dataintelligence-group.com/")
$ws.Hyperlinks.Add($ws.Range("B83"), "https://rpubs.com/dataintelligence/")
$ws.Hyperlinks.Add($ws.Range("B84"), "https://github.com/Sud-Austral/")
$ws.Hyperlinks.Add($ws.Range("B85"), "https://ezexporter.highviewapps.com/exports/export-profile/")

# Hyperlinks.Add() re-styles the target cell with an equivalent but
# distinct "Hyperlink" style record; restore the exact style used by
# every other URL cell in the sheet.
$ws.Range("B72:B85").Style = $ws.Range("B58").Style

Write-Host "Disponibilidad: appended batch rows 72:85, refreshed rows 58:71 timestamp"
